# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 22 de Mayo de 2020 a las 15:05"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1621727
$ws.Range("C4").Value = 825
$ws.Range("E4").Value = 1143106

# Row 11 - Alemania
$ws.Range("B11").Value = 179160
$ws.Range("C11").Value = 139
$ws.Range("E11").Value = 11844

# Row 18 - Arabia Saudita
$ws.Range("B18").Value = 67719
$ws.Range("C18").Value = 2642
$ws.Range("D18").Value = 39003
$ws.Range("E18").Value = 28352
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = 364

# Row 24 - Catar
$ws.Range("B24").Value = 40481
$ws.Range("C24").Value = 1830
$ws.Range("D24").Value = 7893
$ws.Range("E24").Value = 32569
$ws.Range("G24").Value = 2
$ws.Range("H24").Value = 19

# Row 57 - Kazajistan
$ws.Range("D57").Value = 3958
$ws.Range("E57").Value = 3604

# Row 102 - Kenia
$ws.Range("D102").Value = 380
$ws.Range("E102").Value = 731

# Row 125 - San Marino
$ws.Range("B125").Value = 661
$ws.Range("C125").Value = 3
$ws.Range("D125").Value = 254
$ws.Range("E125").Value = 366
